$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new sheet right after "Employee" and make it the active sheet.
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "NegativeLogins"

# Header row
$newSheet.Range("A1").Value = "username"
$newSheet.Range("B1").Value = "password"
$newSheet.Range("C1").Value = "errorMessage"

# Row 2
$newSheet.Range("A2").Value = "Admin"
$newSheet.Range("B2").Value = "admin12"
$newSheet.Range("C2").Value = "Invalid credentials"
$newSheet.Range("D2").NumberFormat = "General"

# Row 3
$newSheet.Range("A3").Value = "Admi"
$newSheet.Range("B3").Value = "admin123"
$newSheet.Range("C3").Value = "Invalid credentials"
$newSheet.Range("D3").NumberFormat = "General"

# Row 4
$newSheet.Range("A4").Value = "Admi"
$newSheet.Range("B4").Value = "admin14"
$newSheet.Range("C4").Value = "Invalid credentials"
$newSheet.Range("D4").NumberFormat = "General"

# Row 5
$newSheet.Range("B5").Value = "admin15"
$newSheet.Range("C5").Value = "Username cannot be empty"
$newSheet.Range("D5").NumberFormat = "General"

# Row 6
$newSheet.Range("A6").Value = "Admin"
$newSheet.Range("C6").Value = "Password cannot be empty"
$newSheet.Range("D6").NumberFormat = "General"

# Row 7
$newSheet.Range("C7").Value = "Username cannot be empty"
$newSheet.Range("D7").NumberFormat = "General"

# Column widths to match the authored layout (closest reachable values —
# the host snaps ColumnWidth to 1/6-character granularity)
$newSheet.Columns.Item(2).ColumnWidth = 8.0
$newSheet.Columns.Item(3).ColumnWidth = 23.5

# Page / print setup to mirror the other sheet in the workbook
$newSheet.PageSetup.PrintHeadings = $false
$newSheet.PageSetup.PrintGridlines = $false
$newSheet.PageSetup.LeftMargin = 50.45669291338584
$newSheet.PageSetup.RightMargin = 50.45669291338584
$newSheet.PageSetup.TopMargin = 54.14173228346456
$newSheet.PageSetup.BottomMargin = 54.14173228346456
$newSheet.PageSetup.HeaderMargin = 21.6
$newSheet.PageSetup.FooterMargin = 21.6
$newSheet.PageSetup.PaperSize = 9
$newSheet.PageSetup.Zoom = 100
$newSheet.PageSetup.FitToPagesWide = 1
$newSheet.PageSetup.FitToPagesTall = 1
$newSheet.PageSetup.Order = 2
$newSheet.PageSetup.Orientation = 1
